# Weekly update: insert two new price rows (new week of data) at the top of
# the data block (rows 26-27), pushing all existing rows from 26..108 down
# to 28..110. This mirrors the source XML diff, which shows every row from
# the old 26..108 reappearing unchanged two rows further down, with the
# dimension growing from A1:R108 to A1:R110.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 26 - this shifts rows
# 26..108 down to 28..110, matching the diff exactly.
$ws.Rows("26:27").Insert()

# Populate the two new rows with the new week's data.
$ws.Range("A26").Value = 2
$ws.Range("B26").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C26").Value = "Coquimbo"
$ws.Range("D26").Value = 44581
$ws.Range("E26").Value = 4
$ws.Range("F26").Value = 100112024
$ws.Range("G26").Value = "Choclo"
$ws.Range("H26").Value = "Choclero"
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 100000
$ws.Range("K26").Value = 170
$ws.Range("L26").Value = 200
$ws.Range("M26").Value = 185
$ws.Range("N26").Value = "$/unidad"
$ws.Range("O26").Value = "Provincia de Limarí"
$ws.Range("P26").Value = 185
$ws.Range("Q26").Value = 1
$ws.Range("R26").Value = "Hortaliza"

$ws.Range("A27").Value = 2
$ws.Range("B27").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C27").Value = "Coquimbo"
$ws.Range("D27").Value = 44581
$ws.Range("E27").Value = 4
$ws.Range("F27").Value = 100112024
$ws.Range("G27").Value = "Choclo"
$ws.Range("H27").Value = "Dulce o Americano"
$ws.Range("I27").Value = "Primera"
$ws.Range("J27").Value = 60000
$ws.Range("K27").Value = 150
$ws.Range("L27").Value = 170
$ws.Range("M27").Value = 160
$ws.Range("N27").Value = "$/unidad"
$ws.Range("O27").Value = "Provincia de Limarí"
$ws.Range("P27").Value = 160
$ws.Range("Q27").Value = 1
$ws.Range("R27").Value = "Hortaliza"
